$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 196, shifting rows 196:228 down to 197:229
$ws.Rows.Item(196).Insert()

# Populate the new row 196 with its own unique data; all "label" columns
# (A,B,C,E,F,G,H,I,N,Q,R) are identical to the surrounding rows.
$ws.Range("A196").Value = 10
$ws.Range("B196").Value = "Vega Modelo de Temuco"
$ws.Range("C196").Value = "La Araucanía"
$ws.Range("D196").Value = 44522
$ws.Range("E196").Value = 9
$ws.Range("F196").Value = 100112044
$ws.Range("G196").Value = "Perejil"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 30
$ws.Range("K196").Value = 5000
$ws.Range("L196").Value = 5000
$ws.Range("M196").Value = 5000
$ws.Range("N196").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O196").Value = "Provincia de Cautín"
$ws.Range("P196").Value = 1667
$ws.Range("Q196").Value = 3
$ws.Range("R196").Value = "Hortaliza"
